# Regenerate the save_data sheet's K column (column G) so it reflects the
# recomputed strike count ("K instead of Strike#") rather than the old
# Strike# value, as part of regenerating std/mean and writing s_vals.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New K values for rows 2-58 (row index => new value), keyed by the sheet's
# data row number so it lines up 1:1 with the regenerated save_data rows.
$kValues = @{
    2  = 0
    3  = 0
    4  = 2
    5  = 5
    6  = 3
    7  = 1
    8  = 0
    9  = 1
    10 = 1
    11 = 2
    12 = 1
    13 = 1
    14 = 1
    15 = 1
    16 = 2
    17 = 1
    18 = 2
    19 = 1
    20 = 0
    21 = 2
    22 = 1
    23 = 1
    24 = 1
    25 = 0
    26 = 2
    27 = 1
    28 = 1
    29 = 1
    30 = 0
    31 = 2
    32 = 0
    33 = 3
    34 = 4
    35 = 1
    36 = 0
    37 = 1
    38 = 2
    39 = 2
    40 = 0
    41 = 0
    42 = 3
    43 = 0
    44 = 1
    45 = 3
    46 = 3
    47 = 2
    48 = 1
    49 = 0
    50 = 2
    51 = 0
    52 = 2
    53 = 0
    54 = 1
    55 = 2
    56 = 0
    57 = 0
    58 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
